# Auto-generated edit script: updates Leve profit/price columns (H-N)
# for the rows identified in the source diff, across 8 worksheets.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Cells.Item(62, 8).Value = 4825
$ws.Cells.Item(62, 9).Value = 4500
$ws.Cells.Item(62, 10).Value = 4933.3335
$ws.Cells.Item(62, 11).Value = 4500
$ws.Cells.Item(62, 12).Value = 4933.3335
$ws.Cells.Item(62, 13).Value = -3876
$ws.Cells.Item(62, 14).Value = -6181.3335
# Row 65
$ws.Cells.Item(65, 8).Value = 4825
$ws.Cells.Item(65, 9).Value = 4500
$ws.Cells.Item(65, 10).Value = 4933.3335
$ws.Cells.Item(65, 11).Value = 22500
$ws.Cells.Item(65, 12).Value = 24666.6675
$ws.Cells.Item(65, 13).Value = -19380
$ws.Cells.Item(65, 14).Value = -30906.6675
# Row 86
$ws.Cells.Item(86, 8).Value = 4382
$ws.Cells.Item(86, 9).Value = 5408.5
$ws.Cells.Item(86, 10).Value = 3635.4546
$ws.Cells.Item(86, 11).Value = 5408.5
$ws.Cells.Item(86, 12).Value = 3635.4546
$ws.Cells.Item(86, 13).Value = -4285.5
$ws.Cells.Item(86, 14).Value = -5881.4546
# Row 89
$ws.Cells.Item(89, 8).Value = 4382
$ws.Cells.Item(89, 9).Value = 5408.5
$ws.Cells.Item(89, 10).Value = 3635.4546
$ws.Cells.Item(89, 11).Value = 27042.5
$ws.Cells.Item(89, 12).Value = 18177.273
$ws.Cells.Item(89, 13).Value = -21426.5
$ws.Cells.Item(89, 14).Value = -29409.273
# Row 106
$ws.Cells.Item(106, 8).Value = 8124.778
$ws.Cells.Item(106, 9).Value = 1562
$ws.Cells.Item(106, 10).Value = 9999.857
$ws.Cells.Item(106, 11).Value = 1562
$ws.Cells.Item(106, 12).Value = 9999.857
$ws.Cells.Item(106, 13).Value = -931
$ws.Cells.Item(106, 14).Value = -11261.857
# Row 107
$ws.Cells.Item(107, 8).Value = 1708.0588
$ws.Cells.Item(107, 9).Value = 1634.2858
$ws.Cells.Item(107, 10).Value = 2052.3333
$ws.Cells.Item(107, 11).Value = 1634.2858
$ws.Cells.Item(107, 12).Value = 2052.3333
$ws.Cells.Item(107, 13).Value = 285.7141999999999
$ws.Cells.Item(107, 14).Value = -5892.3333
# Row 113
$ws.Cells.Item(113, 8).Value = 5002.6
$ws.Cells.Item(113, 9).Value = 4787.25
$ws.Cells.Item(113, 10).Value = 5146.1665
$ws.Cells.Item(113, 11).Value = 4787.25
$ws.Cells.Item(113, 12).Value = 5146.1665
$ws.Cells.Item(113, 13).Value = -1533.25
$ws.Cells.Item(113, 14).Value = -11654.1665
# Row 132
$ws.Cells.Item(132, 8).Value = 3970.5
$ws.Cells.Item(132, 9).Value = 3978.3914
$ws.Cells.Item(132, 10).Value = 3789
$ws.Cells.Item(132, 11).Value = 11935.1742
$ws.Cells.Item(132, 12).Value = 11367
$ws.Cells.Item(132, 13).Value = -9405.174199999999
$ws.Cells.Item(132, 14).Value = -16427

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 6490.344
$ws.Cells.Item(32, 9).Value = 668.4400000000001
$ws.Cells.Item(32, 10).Value = 32953.547
$ws.Cells.Item(32, 11).Value = 668.4400000000001
$ws.Cells.Item(32, 12).Value = 32953.547
$ws.Cells.Item(32, 13).Value = -381.4400000000001
$ws.Cells.Item(32, 14).Value = -33527.547
# Row 45
$ws.Cells.Item(45, 8).Value = 14997.944
$ws.Cells.Item(45, 9).Value = 24740.445
$ws.Cells.Item(45, 10).Value = 5255.4443
$ws.Cells.Item(45, 11).Value = 24740.445
$ws.Cells.Item(45, 12).Value = 5255.4443
$ws.Cells.Item(45, 13).Value = -24363.445
$ws.Cells.Item(45, 14).Value = -6009.4443
# Row 61
$ws.Cells.Item(61, 8).Value = 6370.222
$ws.Cells.Item(61, 9).Value = 6010.9565
$ws.Cells.Item(61, 10).Value = 8436
$ws.Cells.Item(61, 11).Value = 6010.9565
$ws.Cells.Item(61, 12).Value = 8436
$ws.Cells.Item(61, 13).Value = -5798.9565
$ws.Cells.Item(61, 14).Value = -8860
# Row 74
$ws.Cells.Item(74, 8).Value = 3760.625
$ws.Cells.Item(74, 9).Value = 3081.5
$ws.Cells.Item(74, 10).Value = 5798
$ws.Cells.Item(74, 11).Value = 3081.5
$ws.Cells.Item(74, 12).Value = 5798
$ws.Cells.Item(74, 13).Value = -2207.5
$ws.Cells.Item(74, 14).Value = -7546
# Row 77
$ws.Cells.Item(77, 8).Value = 3760.625
$ws.Cells.Item(77, 9).Value = 3081.5
$ws.Cells.Item(77, 10).Value = 5798
$ws.Cells.Item(77, 11).Value = 15407.5
$ws.Cells.Item(77, 12).Value = 28990
$ws.Cells.Item(77, 13).Value = -11039.5
$ws.Cells.Item(77, 14).Value = -37726
# Row 81
$ws.Cells.Item(81, 8).Value = 89999.5
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 89999.5
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 12).Value = 89999.5
$ws.Cells.Item(81, 14).Value = -91995.5
# Row 84
$ws.Cells.Item(84, 8).Value = 89999.5
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 89999.5
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 269998.5
$ws.Cells.Item(84, 14).Value = -279982.5
# Row 132
$ws.Cells.Item(132, 8).Value = 6168.0713
$ws.Cells.Item(132, 9).Value = 6083.154
$ws.Cells.Item(132, 10).Value = 7272
$ws.Cells.Item(132, 11).Value = 18249.462
$ws.Cells.Item(132, 12).Value = 21816
$ws.Cells.Item(132, 13).Value = -15719.462
$ws.Cells.Item(132, 14).Value = -26876
# Row 136
$ws.Cells.Item(136, 8).Value = 6370.222
$ws.Cells.Item(136, 9).Value = 6010.9565
$ws.Cells.Item(136, 10).Value = 8436
$ws.Cells.Item(136, 11).Value = 18032.8695
$ws.Cells.Item(136, 12).Value = 25308
$ws.Cells.Item(136, 13).Value = -15482.8695
$ws.Cells.Item(136, 14).Value = -30408

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 2354.182
$ws.Cells.Item(20, 9).Value = 2273.3157
$ws.Cells.Item(20, 10).Value = 2866.3333
$ws.Cells.Item(20, 11).Value = 2273.3157
$ws.Cells.Item(20, 12).Value = 2866.3333
$ws.Cells.Item(20, 13).Value = -2026.3157
$ws.Cells.Item(20, 14).Value = -3360.3333
# Row 134
$ws.Cells.Item(134, 8).Value = 2693.647
$ws.Cells.Item(134, 9).Value = 2693.647
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 8080.941
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -5545.941
$ws.Cells.Item(134, 14).Value = ""

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 9525
$ws.Cells.Item(31, 9).Value = 12707.417
$ws.Cells.Item(31, 10).Value = 3160.1667
$ws.Cells.Item(31, 11).Value = 12707.417
$ws.Cells.Item(31, 12).Value = 3160.1667
$ws.Cells.Item(31, 13).Value = -12412.417
$ws.Cells.Item(31, 14).Value = -3750.1667
# Row 34
$ws.Cells.Item(34, 8).Value = 9525
$ws.Cells.Item(34, 9).Value = 12707.417
$ws.Cells.Item(34, 10).Value = 3160.1667
$ws.Cells.Item(34, 11).Value = 12707.417
$ws.Cells.Item(34, 12).Value = 3160.1667
$ws.Cells.Item(34, 13).Value = -12505.417
$ws.Cells.Item(34, 14).Value = -3564.1667
# Row 35
$ws.Cells.Item(35, 8).Value = 811.8
$ws.Cells.Item(35, 9).Value = 811.8
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 811.8
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -517.8
$ws.Cells.Item(35, 14).Value = ""
# Row 58
$ws.Cells.Item(58, 8).Value = 1925
$ws.Cells.Item(58, 9).Value = 2000
$ws.Cells.Item(58, 10).Value = 1914.2858
$ws.Cells.Item(58, 11).Value = 2000
$ws.Cells.Item(58, 12).Value = 1914.2858
$ws.Cells.Item(58, 13).Value = -1797
$ws.Cells.Item(58, 14).Value = -2320.2858
# Row 86
$ws.Cells.Item(86, 8).Value = 25648232
$ws.Cells.Item(86, 9).Value = 47624900
$ws.Cells.Item(86, 10).Value = 8787.333000000001
$ws.Cells.Item(86, 11).Value = 47624900
$ws.Cells.Item(86, 12).Value = 8787.333000000001
$ws.Cells.Item(86, 13).Value = -47623777
$ws.Cells.Item(86, 14).Value = -11033.333
# Row 89
$ws.Cells.Item(89, 8).Value = 25648232
$ws.Cells.Item(89, 9).Value = 47624900
$ws.Cells.Item(89, 10).Value = 8787.333000000001
$ws.Cells.Item(89, 11).Value = 238124500
$ws.Cells.Item(89, 12).Value = 43936.665
$ws.Cells.Item(89, 13).Value = -238118884
$ws.Cells.Item(89, 14).Value = -55168.665
# Row 94
$ws.Cells.Item(94, 8).Value = 993.1
$ws.Cells.Item(94, 9).Value = 464.66666
$ws.Cells.Item(94, 10).Value = 1219.5714
$ws.Cells.Item(94, 11).Value = 464.66666
$ws.Cells.Item(94, 12).Value = 1219.5714
$ws.Cells.Item(94, 13).Value = -13.66665999999998
$ws.Cells.Item(94, 14).Value = -2121.5714
# Row 134
$ws.Cells.Item(134, 8).Value = 4965.8823
$ws.Cells.Item(134, 9).Value = 5880.231
$ws.Cells.Item(134, 10).Value = 1994.25
$ws.Cells.Item(134, 11).Value = 17640.693
$ws.Cells.Item(134, 12).Value = 5982.75
$ws.Cells.Item(134, 13).Value = -15105.693
$ws.Cells.Item(134, 14).Value = -11052.75
# Row 136
$ws.Cells.Item(136, 8).Value = 1925
$ws.Cells.Item(136, 9).Value = 2000
$ws.Cells.Item(136, 10).Value = 1914.2858
$ws.Cells.Item(136, 11).Value = 6000
$ws.Cells.Item(136, 12).Value = 5742.857400000001
$ws.Cells.Item(136, 13).Value = -3450
$ws.Cells.Item(136, 14).Value = -10842.8574

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 35
$ws.Cells.Item(35, 8).Value = 500
$ws.Cells.Item(35, 9).Value = 500
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 1500
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -1212
# Row 86
$ws.Cells.Item(86, 8).Value = 730.6
$ws.Cells.Item(86, 9).Value = 998
$ws.Cells.Item(86, 10).Value = 552.3333
$ws.Cells.Item(86, 11).Value = 2994
$ws.Cells.Item(86, 12).Value = 1656.9999
$ws.Cells.Item(86, 13).Value = -1808
$ws.Cells.Item(86, 14).Value = -4028.9999
# Row 89
$ws.Cells.Item(89, 8).Value = 730.6
$ws.Cells.Item(89, 9).Value = 998
$ws.Cells.Item(89, 10).Value = 552.3333
$ws.Cells.Item(89, 11).Value = 8982
$ws.Cells.Item(89, 12).Value = 4970.9997
$ws.Cells.Item(89, 13).Value = -3054
$ws.Cells.Item(89, 14).Value = -16826.9997
# Row 128
$ws.Cells.Item(128, 8).Value = 166015
$ws.Cells.Item(128, 9).Value = 166015
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 498045
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 13).Value = -493065
# Row 140
$ws.Cells.Item(140, 8).Value = 1670
$ws.Cells.Item(140, 9).Value = 1026.2
$ws.Cells.Item(140, 10).Value = 2474.75
$ws.Cells.Item(140, 11).Value = 3078.6
$ws.Cells.Item(140, 12).Value = 7424.25
$ws.Cells.Item(140, 13).Value = 2101.4
$ws.Cells.Item(140, 14).Value = -17784.25
# Row 141
$ws.Cells.Item(141, 8).Value = 5024.8335
$ws.Cells.Item(141, 9).Value = 4030
$ws.Cells.Item(141, 10).Value = 9999
$ws.Cells.Item(141, 11).Value = 12090
$ws.Cells.Item(141, 12).Value = 29997
$ws.Cells.Item(141, 13).Value = -6910
$ws.Cells.Item(141, 14).Value = -40357

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Cells.Item(41, 8).Value = 22010
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 22010
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 22010
$ws.Cells.Item(41, 13).Value = ""
$ws.Cells.Item(41, 14).Value = -22720
# Row 102
$ws.Cells.Item(102, 8).Value = 1630
$ws.Cells.Item(102, 9).Value = 1657.2858
$ws.Cells.Item(102, 10).Value = 1566.3334
$ws.Cells.Item(102, 11).Value = 1657.2858
$ws.Cells.Item(102, 12).Value = 1566.3334
$ws.Cells.Item(102, 13).Value = -35.28580000000011
$ws.Cells.Item(102, 14).Value = -4810.3334
# Row 132
$ws.Cells.Item(132, 8).Value = 5296
$ws.Cells.Item(132, 9).Value = 4823.143
$ws.Cells.Item(132, 10).Value = 6123.5
$ws.Cells.Item(132, 11).Value = 14469.429
$ws.Cells.Item(132, 12).Value = 18370.5
$ws.Cells.Item(132, 13).Value = -11939.429
$ws.Cells.Item(132, 14).Value = -23430.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 15706.956
$ws.Cells.Item(7, 9).Value = 13360.048
$ws.Cells.Item(7, 10).Value = 40349.5
$ws.Cells.Item(7, 11).Value = 13360.048
$ws.Cells.Item(7, 12).Value = 40349.5
$ws.Cells.Item(7, 13).Value = -13248.048
$ws.Cells.Item(7, 14).Value = -40573.5
# Row 22
$ws.Cells.Item(22, 8).Value = 2202.0312
$ws.Cells.Item(22, 9).Value = 2022.75
$ws.Cells.Item(22, 10).Value = 2381.3125
$ws.Cells.Item(22, 11).Value = 2022.75
$ws.Cells.Item(22, 12).Value = 2381.3125
$ws.Cells.Item(22, 13).Value = -1727.75
$ws.Cells.Item(22, 14).Value = -2971.3125
# Row 27
$ws.Cells.Item(27, 8).Value = 2202.0312
$ws.Cells.Item(27, 9).Value = 2022.75
$ws.Cells.Item(27, 10).Value = 2381.3125
$ws.Cells.Item(27, 11).Value = 2022.75
$ws.Cells.Item(27, 12).Value = 2381.3125
$ws.Cells.Item(27, 13).Value = -1915.75
$ws.Cells.Item(27, 14).Value = -2595.3125
# Row 46
$ws.Cells.Item(46, 8).Value = 4895.909
$ws.Cells.Item(46, 9).Value = 5693.3335
$ws.Cells.Item(46, 10).Value = 3939
$ws.Cells.Item(46, 11).Value = 5693.3335
$ws.Cells.Item(46, 12).Value = 3939
$ws.Cells.Item(46, 13).Value = -5505.3335
$ws.Cells.Item(46, 14).Value = -4315
# Row 100
$ws.Cells.Item(100, 8).Value = 3638.1765
$ws.Cells.Item(100, 9).Value = 1984.9
$ws.Cells.Item(100, 10).Value = 6000
$ws.Cells.Item(100, 11).Value = 1984.9
$ws.Cells.Item(100, 12).Value = 6000
$ws.Cells.Item(100, 13).Value = -1443.9
$ws.Cells.Item(100, 14).Value = -7082
# Row 126
$ws.Cells.Item(126, 8).Value = 15706.956
$ws.Cells.Item(126, 9).Value = 13360.048
$ws.Cells.Item(126, 10).Value = 40349.5
$ws.Cells.Item(126, 11).Value = 40080.144
$ws.Cells.Item(126, 12).Value = 121048.5
$ws.Cells.Item(126, 13).Value = -37610.144
$ws.Cells.Item(126, 14).Value = -125988.5
# Row 136
$ws.Cells.Item(136, 8).Value = 2280.125
$ws.Cells.Item(136, 9).Value = 2040.3334
$ws.Cells.Item(136, 10).Value = 2999.5
$ws.Cells.Item(136, 11).Value = 6121.0002
$ws.Cells.Item(136, 12).Value = 8998.5
$ws.Cells.Item(136, 13).Value = -3571.0002
$ws.Cells.Item(136, 14).Value = -14098.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = ""
# Row 84
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = ""
# Row 86
$ws.Cells.Item(86, 8).Value = 47499.5
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 47499.5
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 47499.5
$ws.Cells.Item(86, 14).Value = -49745.5
# Row 89
$ws.Cells.Item(89, 8).Value = 47499.5
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 47499.5
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 237497.5
$ws.Cells.Item(89, 14).Value = -248729.5
# Row 126
$ws.Cells.Item(126, 8).Value = 1901.05
$ws.Cells.Item(126, 9).Value = 1218.2858
$ws.Cells.Item(126, 10).Value = 2268.6924
$ws.Cells.Item(126, 11).Value = 3654.8574
$ws.Cells.Item(126, 12).Value = 6806.0772
$ws.Cells.Item(126, 13).Value = -1184.8574
$ws.Cells.Item(126, 14).Value = -11746.0772
# Row 132
$ws.Cells.Item(132, 8).Value = 3493.7856
$ws.Cells.Item(132, 9).Value = 3400.25
$ws.Cells.Item(132, 10).Value = 4055
$ws.Cells.Item(132, 11).Value = 10200.75
$ws.Cells.Item(132, 12).Value = 12165
$ws.Cells.Item(132, 13).Value = -7670.75
$ws.Cells.Item(132, 14).Value = -17225

Write-Output "Applied 337 cell updates and 5 cell clears."
